$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell is set to text (matching the workbook's original inline-string
# cell type) instead of being auto-coerced to a number by COM when the new value
# happens to look numeric.

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "hh1"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "root.pop_gridConnections[0]"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "0.3"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "0.602017313814645"
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "-554.9071135646353"
$ws.Range("R2").NumberFormat = "@"
$ws.Range("R2").Value = "-207.2880356929771"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "hh2"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "root.pop_gridConnections[1]"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "0.3"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "0.602017313814645"
$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value = "-227.6011135646638"
$ws.Range("R3").NumberFormat = "@"
$ws.Range("R3").Value = "-207.2880356929771"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "hh3"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "root.pop_gridConnections[2]"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "root.pop_energySuppliers[1]"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "0.3"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "0.5776933819433462"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "0.0"
$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value = "-163.71911356467206"
$ws.Range("Q4").NumberFormat = "@"
$ws.Range("Q4").Value = "-397.82552304712624"
$ws.Range("R4").NumberFormat = "@"
$ws.Range("R4").Value = "0.0"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "hh4"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "hol1"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "root.pop_energyHolons[0]( p_actorID = hol1, p_actorType = holon, p_parentActorID = sup2 )"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "root.pop_gridConnections[3]"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "0.3"
$ws.Range("J5").NumberFormat = "@"
$ws.Range("J5").Value = "0.602017313814645"
$ws.Range("P5").NumberFormat = "@"
$ws.Range("P5").Value = "-797.4991135646349"
$ws.Range("R5").NumberFormat = "@"
$ws.Range("R5").Value = "-207.2880356929771"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "hh5"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "sup1"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "root.pop_energySuppliers[0]"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "root.pop_gridConnections[4]"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "null"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "0.3"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "0.0"
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = "0.602017313814645"
$ws.Range("P6").NumberFormat = "@"
$ws.Range("P6").Value = "-186.1471135646692"
$ws.Range("Q6").NumberFormat = "@"
$ws.Range("Q6").Value = "0.0"
$ws.Range("R6").NumberFormat = "@"
$ws.Range("R6").Value = "-207.2880356929771"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "hh6"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "root.pop_gridConnections[5]"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "0.3"
$ws.Range("J7").NumberFormat = "@"
$ws.Range("J7").Value = "0.602017313814645"
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "-1164.831113564687"
$ws.Range("R7").NumberFormat = "@"
$ws.Range("R7").Value = "-207.2880356929771"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "hh7"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "root.pop_gridConnections[6]"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "0.3"
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = "0.602017313814645"
$ws.Range("P8").NumberFormat = "@"
$ws.Range("P8").Value = "-475.6111135646337"
$ws.Range("R8").NumberFormat = "@"
$ws.Range("R8").Value = "-207.2880356929771"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "hh8"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "household"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "root.pop_gridConnections[7]"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "0.3"
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = "0.602017313814645"
$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "-450.15911356463477"
$ws.Range("R9").NumberFormat = "@"
$ws.Range("R9").Value = "-207.2880356929771"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "hh9"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "root.pop_gridConnections[8]"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "0.3"
$ws.Range("J10").NumberFormat = "@"
$ws.Range("J10").Value = "0.602017313814645"
$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value = "-267.0391135646595"
$ws.Range("R10").NumberFormat = "@"
$ws.Range("R10").Value = "-207.2880356929771"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "hh10"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "hol1"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "root.pop_energyHolons[0]( p_actorID = hol1, p_actorType = holon, p_parentActorID = sup2 )"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "root.pop_gridConnections[9]"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "root.pop_energySuppliers[1]"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "0.3"
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "0.5776933819433462"
$ws.Range("J11").NumberFormat = "@"
$ws.Range("J11").Value = "0.0"
$ws.Range("P11").NumberFormat = "@"
$ws.Range("P11").Value = "-377.75111356463697"
$ws.Range("Q11").NumberFormat = "@"
$ws.Range("Q11").Value = "-397.82552304712624"
$ws.Range("R11").NumberFormat = "@"
$ws.Range("R11").Value = "0.0"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "hh11"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "household"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "root.pop_gridConnections[10]"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "0.3"
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = "0.5776933819433462"
$ws.Range("P12").NumberFormat = "@"
$ws.Range("P12").Value = "-135.32711356467314"
$ws.Range("Q12").NumberFormat = "@"
$ws.Range("Q12").Value = "-397.82552304712624"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "0.3"
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "0.5776933819433462"
$ws.Range("P13").NumberFormat = "@"
$ws.Range("P13").Value = "-374.1811135646364"
$ws.Range("Q13").NumberFormat = "@"
$ws.Range("Q13").Value = "-397.82552304712624"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "hh13"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "root.pop_gridConnections[12]"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "null"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "0.3"
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "0.0"
$ws.Range("J14").NumberFormat = "@"
$ws.Range("J14").Value = "0.602017313814645"
$ws.Range("P14").NumberFormat = "@"
$ws.Range("P14").Value = "-218.31911356466543"
$ws.Range("Q14").NumberFormat = "@"
$ws.Range("Q14").Value = "0.0"
$ws.Range("R14").NumberFormat = "@"
$ws.Range("R14").Value = "-207.2880356929771"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "hh14"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "root.pop_gridConnections[13]"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "0.3"
$ws.Range("J15").NumberFormat = "@"
$ws.Range("J15").Value = "0.602017313814645"
$ws.Range("P15").NumberFormat = "@"
$ws.Range("P15").Value = "-500.3491135646339"
$ws.Range("R15").NumberFormat = "@"
$ws.Range("R15").Value = "-207.2880356929771"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "hh15"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "sup1"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "root.pop_energySuppliers[0]"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "root.pop_gridConnections[14]"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "0.3"
$ws.Range("J16").NumberFormat = "@"
$ws.Range("J16").Value = "0.602017313814645"
$ws.Range("P16").NumberFormat = "@"
$ws.Range("P16").Value = "-431.3431135646344"
$ws.Range("R16").NumberFormat = "@"
$ws.Range("R16").Value = "-207.2880356929771"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "hh16"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "root.pop_gridConnections[15]"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "0.3"
$ws.Range("J17").NumberFormat = "@"
$ws.Range("J17").Value = "0.602017313814645"
$ws.Range("P17").NumberFormat = "@"
$ws.Range("P17").Value = "-296.9431135646496"
$ws.Range("R17").NumberFormat = "@"
$ws.Range("R17").Value = "-207.2880356929771"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "hh17"
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "root.pop_gridConnections[16]"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "0.3"
$ws.Range("J18").NumberFormat = "@"
$ws.Range("J18").Value = "0.602017313814645"
$ws.Range("P18").NumberFormat = "@"
$ws.Range("P18").Value = "-280.1011135646563"
$ws.Range("R18").NumberFormat = "@"
$ws.Range("R18").Value = "-207.2880356929771"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "hh18"
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "root.pop_gridConnections[17]"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "0.3"
$ws.Range("J19").NumberFormat = "@"
$ws.Range("J19").Value = "0.602017313814645"
$ws.Range("P19").NumberFormat = "@"
$ws.Range("P19").Value = "-440.6671135646345"
$ws.Range("R19").NumberFormat = "@"
$ws.Range("R19").Value = "-207.2880356929771"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "com1"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "commercial"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "hol1"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "root.pop_energyHolons[0]( p_actorID = hol1, p_actorType = holon, p_parentActorID = sup2 )"
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "root.pop_gridConnections[18]"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "root.pop_energySuppliers[1]"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "0.0"
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = "0.0"
$ws.Range("P20").NumberFormat = "@"
$ws.Range("P20").Value = "0.0"
$ws.Range("R20").NumberFormat = "@"
$ws.Range("R20").Value = "0.0"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "com2"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "commercial"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "sup1"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "root.pop_energySuppliers[0]"
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "root.pop_gridConnections[19]"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "null"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "0.0"
$ws.Range("I21").NumberFormat = "@"
$ws.Range("I21").Value = "0.0"
$ws.Range("P21").NumberFormat = "@"
$ws.Range("P21").Value = "0.0"
$ws.Range("Q21").NumberFormat = "@"
$ws.Range("Q21").Value = "0.0"
